# Coureurs sheet: add the remaining 2023 F1 calendar races as new columns
# (I:AB), each initialised to 0 race points, and reset the old G column
# (which had accidentally been filled with a copy of the F/Bahrein value)
# back to 0 now that it represents a still-unplayed race.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coureurs")

# Race names (countries/venues) for columns I..AB (2023 calendar, rounds 4-23)
$raceNames = @(
    "Azerbeidzjan",
    "Miami",
    "Emilia-Romagna",
    "Monaco",
    "Spanje",
    "Canada",
    "Oostenrijk",
    "Groot-Brittanië",
    "Hongarije",
    "België",
    "Nederland",
    "Italië",
    "Singapore",
    "Japan",
    "Qatar",
    "Verenigde Staten",
    "Mexico",
    "Brazilië",
    "Las Vegas",
    "Abu Dhabi"
)

$firstNewCol = 9   # column I
$lastRow = 21

# Header row: write the new race name headers into I1:AB1
for ($i = 0; $i -lt $raceNames.Length; $i++) {
    $ws.Cells.Item(1, $firstNewCol + $i).Value = $raceNames[$i]
}

# Reset column G (Saudi-Arabië results) back to 0 for every driver row,
# and zero-fill all of the brand-new race columns (I:AB) for every row.
for ($row = 2; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 7).Value = 0   # column G

    for ($i = 0; $i -lt $raceNames.Length; $i++) {
        $ws.Cells.Item($row, $firstNewCol + $i).Value = 0
    }
}

# Update the view: scroll so column J is leftmost and select S13
$win = $excel.ActiveWindow
$win.ScrollColumn = 10
$win.ScrollRow = 1
$ws.Range("S13").Select()
